$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '70.961.92'
$ws.Range("E2").Value = '  +0.84%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.590.47'
$ws.Range("E3").Value = '  +0.27%  '

$ws.Range("E4").Value = '  -0.33%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '586.92'
$ws.Range("E5").Value = '  +1.60%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '185.43'
$ws.Range("E6").Value = '  -1.06%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.577.15'
$ws.Range("E7").Value = '  +0.55%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.623'
$ws.Range("E8").Value = '  +0.19%  '

$ws.Range("E9").Value = '  -0.24%  '

$ws.Range("E10").Value = '  +16.46%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.650'
$ws.Range("E11").Value = '  +0.23%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.26'
$ws.Range("E12").Value = '  +0.01%  '

$ws.Range("E13").Value = '  +5.55%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.53'
$ws.Range("E14").Value = '  -0.63%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.159.17'
$ws.Range("E15").Value = '  -1.01%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '70.892.06'
$ws.Range("E16").Value = '  +0.56%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '19.32'
$ws.Range("E17").Value = '  -1.09%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.558.31'
$ws.Range("E18").Value = '  -1.20%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.41'
$ws.Range("E19").Value = '  -0.62%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '568.62'
$ws.Range("E20").Value = '  +13.46%  '

$ws.Range("E21").Value = '  -0.04%  '

$ws.Range("E22").Value = '  -3.02%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '17.65'
$ws.Range("E23").Value = '  -9.09%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.62'
$ws.Range("E24").Value = '  +5.78%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.01'
$ws.Range("E25").Value = '  +1.90%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '94.85'
$ws.Range("E26").Value = '  -0.35%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.22'
$ws.Range("E27").Value = '  -2.59%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.93'
$ws.Range("E28").Value = '  -1.20%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.11'
$ws.Range("E29").Value = '  -2.27%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '32.38'
$ws.Range("E30").Value = '  +2.13%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.26'
$ws.Range("E31").Value = '  -5.27%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.31'
$ws.Range("E32").Value = '  -1.80%  '

$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.115'
$ws.Range("E33").Value = '  -1.56%  '

$ws.Range("B34").Value = 'OKB'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '64.10'
$ws.Range("E34").Value = '  -3.09%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.30'
$ws.Range("E35").Value = '  +1.44%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '547.92'
$ws.Range("E36").Value = '  -4.46%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.414'
$ws.Range("E37").Value = '  +0.86%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0₃0807'
$ws.Range("E38").Value = '  +2.07%  '

$ws.Range("E39").Value = '  +0.18%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '37.58'
$ws.Range("E40").Value = '  -3.03%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.479.02'
$ws.Range("E41").Value = '  +8.82%  '

$ws.Range("B42").Value = 'dogwifhat'
$ws.Range("C42").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.17'
$ws.Range("E42").Value = '  -3.02%  '

$ws.Range("B43").Value = 'Kaspa'
$ws.Range("C43").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.136'
$ws.Range("E43").Value = '  -0.39%  '

$ws.Range("B44").Value = 'Stacks'
$ws.Range("C44").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.41'
$ws.Range("E44").Value = '  -0.88%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.52'
$ws.Range("E45").Value = '  -1.07%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.97'
$ws.Range("E46").Value = '  -2.28%  '

$ws.Range("E47").Value = '  -1.95%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.38'
$ws.Range("E48").Value = '  -1.28%  '

$ws.Range("E49").Value = '  +1.63%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.998'
$ws.Range("E50").Value = '  -0.40%  '

$ws.Range("E51").Value = '  -6.01%  '
